$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.995.60'
$ws.Range("E2").Value = '  -2.29%  '
$ws.Range("D3").Value = '3.491.48'
$ws.Range("E3").Value = '  +1.39%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.49'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.99%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.17'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.74%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  -0.97%  '
$ws.Range("D9").Value = '3.488.20'
$ws.Range("E9").Value = '  +1.38%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.131'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.24%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.86'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.41%  '
$ws.Range("E12").Value = '  -3.64%  '
$ws.Range("D13").Value = '4.097.80'
$ws.Range("E13").Value = '  +1.31%  '
$ws.Range("E14").Value = '  +0.90%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '29.96'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.88%  '
$ws.Range("D16").Value = '66.100.46'
$ws.Range("E16").Value = '  -2.11%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000171'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.54%  '
$ws.Range("D18").Value = '3.490.29'
$ws.Range("E18").Value = '  +1.23%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.93'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.90'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.04%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '366.84'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.29%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.75'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.70%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.37%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.27'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.42%  '
$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000126'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.45%  '
$ws.Range("B26").Value = 'Polygon'
$ws.Range("C26").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.536'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.61'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.21%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.180'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.96%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '24.05'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.19%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.77'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.21%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.99'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.51%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.05%  '
$ws.Range("B34").Value = 'Fetch.AI'
$ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.29'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.84%  '
$ws.Range("B35").Value = 'Aptos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.13'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.02%  '
$ws.Range("E36").Value = '  -0.97%  '
$ws.Range("B37").Value = 'EnergySwap'
$ws.Range("C37").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '29.61'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +14.92%  '
$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '159.21'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.64%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.890'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.21%  '
$ws.Range("D40").Value = '2.821.27'
$ws.Range("E40").Value = '  +4.79%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.76'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.03%  '
$ws.Range("E42").Value = '  -5.75%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.45'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.28%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.45'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.46%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0682'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.92%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.15'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.41%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.15'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.65%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0288'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.32%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '317.65'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.47%  '
$ws.Range("E50").Value = '  -2.27%  '
$ws.Range("E51").Value = '  -2.37%  '
